$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "In Translation"
#    (Overview!E2:F2, zh-cn!C2, de-de!C2 all shared this string)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2) Narrow the "status" columns (previously auto-sized wider than needed).
#    ColumnWidth is in characters; the values below reproduce the narrower
#    width from the updated report as closely as this host's column-width
#    model allows.
# ---------------------------------------------------------------------------
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
